$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I (2000) ---
$ws.Range("I12").Value = 4188377156.3100057
$ws.Range("I13").Value = 1012006300.0300001
$ws.Range("I14").Value = -44319159.290000051
$ws.Range("I16").Value = -162861893.56999999
$ws.Range("I18").Formula = "=SUM(I12:I17)"
$ws.Range("I19").Value = -1160500000.0000002
$ws.Range("I21").Formula = "=SUM(I18:I20)"
$ws.Range("I26").Value = 1010658958.9880759

# --- Column J (1999) ---
$ws.Range("J12").Value = 3588029419
$ws.Range("J13").Value = 956934340.60000002
$ws.Range("J14").Value = 146268235.09999999
$ws.Range("J16").Value = 193292161.30000001
$ws.Range("J19").Value = 1105900000
$ws.Range("J22").Value = 57815625
$ws.Range("J26").Value = 1031977291

$wb.Application.Calculate()
